# Auto-generated cell updates for D (Price) and E (Volume(1h)) columns
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'291.03"
$ws.Range("E2").Value = "'-3.28%"
$ws.Range("D3").Value = "'30.71"
$ws.Range("E3").Value = "'-4.78%"
$ws.Range("D4").Value = "'4.948"
$ws.Range("E4").Value = "'-0.09%"
$ws.Range("D5").Value = "'0.07220"
$ws.Range("E5").Value = "'-5.57%"
$ws.Range("D6").Value = "'1.840"
$ws.Range("E6").Value = "'-4.20%"
$ws.Range("D7").Value = "'7.696"
$ws.Range("E7").Value = "'-1.79%"
$ws.Range("D8").Value = "'3.770"
$ws.Range("E8").Value = "'-0.68%"
$ws.Range("D9").Value = "'0.8970"
$ws.Range("E9").Value = "'-2.20%"
$ws.Range("D10").Value = "'0.1656"
$ws.Range("E10").Value = "'-5.17%"
$ws.Range("D11").Value = "'0.07706"
$ws.Range("E11").Value = "'-0.52%"
$ws.Range("D12").Value = "'0.08058"
$ws.Range("E12").Value = "'-5.61%"
$ws.Range("D13").Value = "'0.03044"
$ws.Range("E13").Value = "'-3.97%"
$ws.Range("D14").Value = "'0.1000"
$ws.Range("E14").Value = "'0.19%"
$ws.Range("D15").Value = "'0.001501"
$ws.Range("E15").Value = "'-1.14%"
$ws.Range("D16").Value = "'0.005759"
$ws.Range("E16").Value = "'-2.64%"
$ws.Range("D18").Value = "'3.470"
$ws.Range("E18").Value = "'0.13%"
$ws.Range("E19").Value = "'-3.27%"
$ws.Range("D20").Value = "'0.3318"
$ws.Range("E20").Value = "'-0.95%"
$ws.Range("D21").Value = "'0.1298"
$ws.Range("E21").Value = "'-2.17%"
$ws.Range("D22").Value = "'4.053"
$ws.Range("E22").Value = "'-5.53%"
$ws.Range("E23").Value = "'19.93%"
$ws.Range("D24").Value = "'0.04515"
$ws.Range("E24").Value = "'-0.02%"
$ws.Range("D25").Value = "'0.001216"
$ws.Range("E25").Value = "'-0.52%"
$ws.Range("D26").Value = "'0.004659"
$ws.Range("E26").Value = "'6.22%"
$ws.Range("E27").Value = "'-0.03%"
$ws.Range("D39").Value = "'0.01577"
$ws.Range("E39").Value = "'-7.00%"
$ws.Range("D40").Value = "'0.04394"
$ws.Range("D41").Value = "'0.007267"
$ws.Range("E41").Value = "'-2.58%"
$ws.Range("D42").Value = "'0.01004"
$ws.Range("D43").Value = "'0.1306"
$ws.Range("E43").Value = "'-3.20%"
$ws.Range("E44").Value = "'-11.61%"
$ws.Range("D45").Value = "'0.009196"
$ws.Range("E45").Value = "'-12.25%"
$ws.Range("D46").Value = "'0.00005959"
$ws.Range("E46").Value = "'-4.54%"
$ws.Range("E47").Value = "'-0.07%"
$ws.Range("D48").Value = "'2.310"
$ws.Range("E48").Value = "'180.60%"
$ws.Range("E50").Value = "'-0.07%"
$ws.Range("E51").Value = "'-0.07%"
